# Re-saves the "Logistics Operations Data" workbook the way the target
# revision does: renamed sheets, simplified/"de-styled" cell formatting
# (the blanket centred style that used to cover every column is gone,
# only the 2-decimal number format survives on the numeric tables), and
# refreshed sheet selections/active-tab state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the sheets.
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "Return demand "
$wb.Worksheets.Item(2).Name = "Transportation cost btween c,f "
$wb.Worksheets.Item(3).Name = "Transportation cost btween m,f"
$wb.Worksheets.Item(4).Name = "Installing cost of an equipment"

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# 2. Strip the old blanket "centered" styling from every sheet (this
#    also drops the whole-column <cols style="1"/> declarations that
#    only existed to carry that default alignment).
# ---------------------------------------------------------------------
$ws1.Cells.ClearFormats()
$ws2.Cells.ClearFormats()
$ws3.Cells.ClearFormats()
$ws4.Cells.ClearFormats()

# ---------------------------------------------------------------------
# 3. Re-apply the 2-decimal numeric display format to the two
#    "distance" tables - that's the only formatting the new revision
#    keeps (minus the centred alignment it used to carry).
# ---------------------------------------------------------------------
$ws2.Range("B2:I21").NumberFormat = "0.00"
$ws3.Range("B2:D9").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 4. Refresh each sheet's selection / scroll position. Order matters:
#    whichever sheet we touch last becomes the active ("tabSelected")
#    one, so sheet 1 ("Return demand ") is selected last to match the
#    target workbook state.
# ---------------------------------------------------------------------
$ws2.Range("D24").Select()
$ws3.Range("G31").Select()
$ws4.Range("H28").Select()
$ws1.Range("C31").Select()
